# Adds the "ODI Bowling Extra" worksheet (sheetId 5) after "ODI Batting Extra",
# with MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL columns.

$wb = $excel.ActiveWorkbook

# New sheet goes after the current last sheet ("ODI Batting Extra").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# Force the whole used range to be plain text so numeric-looking values
# ("4423", "0", "20.00%", ...) are stored as text, matching the source data.
$ws.Range("A1:C21").NumberFormat = "@"

# Header row (bold, centered, top-aligned, thin border all round).
$ws.Cells.Item(1,1).Value = "MATCH_CODE"
$ws.Cells.Item(1,2).Value = "MAIDEN_OVERS"
$ws.Cells.Item(1,3).Value = "PERCENT_WICKETS_OF_ALL"
$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Data rows.
$ws.Cells.Item(2,1).Value = "4423"
$ws.Cells.Item(2,2).Value = "0"
$ws.Cells.Item(3,1).Value = "4429"
$ws.Cells.Item(3,2).Value = "0"
$ws.Cells.Item(4,1).Value = "4430"
$ws.Cells.Item(4,2).Value = "1"
$ws.Cells.Item(4,3).Value = "20.00%"
$ws.Cells.Item(5,1).Value = "4431"
$ws.Cells.Item(5,2).Value = "0"
$ws.Cells.Item(5,3).Value = "30.00%"
$ws.Cells.Item(6,1).Value = "4435"
$ws.Cells.Item(7,1).Value = "4436"
$ws.Cells.Item(8,1).Value = "4483"
$ws.Cells.Item(8,2).Value = "1"
$ws.Cells.Item(8,3).Value = "50.00%"
$ws.Cells.Item(9,1).Value = "4484"
$ws.Cells.Item(9,2).Value = "1"
$ws.Cells.Item(9,3).Value = "30.00%"
$ws.Cells.Item(10,1).Value = "4486"
$ws.Cells.Item(11,1).Value = "4644"
$ws.Cells.Item(11,2).Value = "1"
$ws.Cells.Item(11,3).Value = "10.00%"
$ws.Cells.Item(12,1).Value = "4645"
$ws.Cells.Item(13,1).Value = "4646"
$ws.Cells.Item(13,2).Value = "0"
$ws.Cells.Item(13,3).Value = "10.00%"
$ws.Cells.Item(14,1).Value = "4647"
$ws.Cells.Item(15,1).Value = "4648"
$ws.Cells.Item(15,2).Value = "0"
$ws.Cells.Item(15,3).Value = "20.00%"
$ws.Cells.Item(16,1).Value = "4649"
$ws.Cells.Item(16,2).Value = "0"
$ws.Cells.Item(16,3).Value = "30.00%"
$ws.Cells.Item(17,1).Value = "4660"
$ws.Cells.Item(18,1).Value = "4663"
$ws.Cells.Item(19,1).Value = "4725"
$ws.Cells.Item(19,2).Value = "0"
$ws.Cells.Item(19,3).Value = "30.00%"
$ws.Cells.Item(20,1).Value = "4728"
$ws.Cells.Item(20,2).Value = "1"
$ws.Cells.Item(20,3).Value = "50.00%"
$ws.Cells.Item(21,1).Value = "4732"
$ws.Cells.Item(21,2).Value = "0"

$ws.Range("A1").Select()
